$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above row 363 (shifts rows 363:410 down to 364:411).
$ws.Rows.Item(363).Insert()

# Seed the new row with the same record template as the row that just moved
# down to 364 (same market/category/quality/etc.), then overwrite the date.
$src = $ws.Range("A364:R364")
$dst = $ws.Range("A363:R363")
$src.Copy($dst)

# New record's date (serial 44984 = 2023-02-27).
$ws.Cells.Item(363, 4).Value = 44984
